# Actualización automática 2025-05-30 16:20:08
# Refresh the per-salesperson/per-client sales-by-category figures on the
# "VENTAS POR GRUPO" sheet (columns C:N, rows 3-50) plus the "<n> de 50"
# non-zero counters in the totals row (row 52).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 3
$ws.Range("K3").Value = 0

# Row 4
$ws.Range("C4").Value = 4655.24
$ws.Range("D4").Value = 1221.12
$ws.Range("E4").Value = 69.45
$ws.Range("F4").Value = 52.25
$ws.Range("K4").Value = 855.36
$ws.Range("L4").Value = 3689.36
$ws.Range("M4").Value = 497.37

# Row 5
$ws.Range("C5").Value = 513.22
$ws.Range("D5").Value = 0
$ws.Range("H5").Value = 290.56
$ws.Range("K5").Value = 2317.94
$ws.Range("L5").Value = 5502.04

# Row 6
$ws.Range("C6").Value = 1026.43
$ws.Range("L6").Value = 17.85

# Row 7
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0

# Row 10
$ws.Range("N10").Value = 2645.91

# Row 13
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 1520.92

# Row 14
$ws.Range("D14").Value = 3870.71
$ws.Range("L14").Value = 2913.81
$ws.Range("N14").Value = 136.86

# Row 15
$ws.Range("D15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -1609.23

# Row 16
$ws.Range("C16").Value = 4582.65
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 69.45
$ws.Range("F16").Value = 52.25
$ws.Range("K16").Value = 855.36
$ws.Range("L16").Value = 1845.15

# Row 18
$ws.Range("L18").Value = -20.74

# Row 19
$ws.Range("K19").Value = 0

# Row 24
$ws.Range("C24").Value = 518.4
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 69.45
$ws.Range("F24").Value = 52.25
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 76.14
$ws.Range("K24").Value = 2634.2
$ws.Range("L24").Value = 13555.51

# Row 27
$ws.Range("L27").Value = 171.19

# Row 28
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 69.45
$ws.Range("F28").Value = 9.779999999999999
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 2041.93
$ws.Range("N28").Value = 0

# Row 30
$ws.Range("E30").Value = 433.8

# Row 31
$ws.Range("D31").Value = 1173.5
$ws.Range("L31").Value = 1038.12

# Row 33
$ws.Range("C33").Value = 2550.53
$ws.Range("D33").Value = 7193.38
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1448.45
$ws.Range("L33").Value = 1816.45

# Row 34
$ws.Range("D34").Value = 1704.96
$ws.Range("E34").Value = 69.45
$ws.Range("F34").Value = 52.25
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2037.1
$ws.Range("L34").Value = 1664.58

# Row 35
$ws.Range("I35").Value = 158.54

# Row 38
$ws.Range("E38").Value = 69.45
$ws.Range("F38").Value = 52.25
$ws.Range("J38").Value = 548.21
$ws.Range("K38").Value = 2888.67
$ws.Range("L38").Value = 11.68

# Row 39
$ws.Range("J39").Value = -60.91
$ws.Range("L39").Value = 0

# Row 40
$ws.Range("D40").Value = 1208.82
$ws.Range("E40").Value = 69.45
$ws.Range("F40").Value = 52.25
$ws.Range("L40").Value = 1087.75

# Row 44
$ws.Range("L44").Value = 0

# Row 45
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 1179.46

# Row 46
$ws.Range("C46").Value = 1016.06
$ws.Range("D46").Value = 814.08

# Row 49
$ws.Range("I49").Value = 0

# Row 50
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 489.11
$ws.Range("N50").Value = 0

# Row 52 - "<n> de 50" non-zero counters, recomputed for the data above
$ws.Range("C52").Value = "7 de 50"
$ws.Range("D52").Value = "7 de 50"
$ws.Range("E52").Value = "8 de 50"
$ws.Range("F52").Value = "7 de 50"
$ws.Range("G52").Value = "0 de 50"
$ws.Range("H52").Value = "1 de 50"
$ws.Range("I52").Value = "1 de 50"
$ws.Range("J52").Value = "2 de 50"
$ws.Range("K52").Value = "7 de 50"
$ws.Range("L52").Value = "16 de 50"
